$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.937.91"
$ws.Range("E2").Value = "  -0.38%  "

# Row 3
$ws.Range("D3").Value = "1.641.95"
$ws.Range("E3").Value = "  -0.55%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -1.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5061"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.75%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.11%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2577"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.52%  "

# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06421"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.59"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.20%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07751"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.15%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.268"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.15%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.866.34"
$ws.Range("E13").Value = "  -0.65%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.630.11"
$ws.Range("E14").Value = "  -1.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5470"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.16%  "

# Row 16
$ws.Range("D16").Value = "0.0₅7943"
$ws.Range("E16").Value = "  -0.41%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.48"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.19%  "

# Row 18
$ws.Range("D18").Value = "25.958.70"
$ws.Range("E18").Value = "  -0.46%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.82%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.379"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.895"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.39%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.979"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.84%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.860"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.21%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.04"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.85%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1138"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.45%  "

# Row 28
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.802"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.74%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.66"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.79%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.244"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.06%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04911"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.39%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.272"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.211"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.57%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.547"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.16%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.365"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.06%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8942"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.31%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.630"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.16%  "

# Row 38
$ws.Range("D38").Value = "1.153.12"
$ws.Range("E38").Value = "  -0.84%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5597"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.43%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01566"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.10%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.713"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.69%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.76"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8050"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.34%  "

# Row 45
$ws.Range("D45").Value = "1.776.87"
$ws.Range("E45").Value = "  -0.70%  "

# Row 46
$ws.Range("D46").Value = "0.0₈118"
$ws.Range("E46").Value = "  +5.21%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4529"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.77%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.74%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.70"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.79%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05052"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.64%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("D51").ClearFormats()

